{"js": "// The document contains a single table of arithmetic expressions\n// (20 rows x 5 columns = 100 cells, one \"<w:t>\" run per cell). The\n// commit replaces each cell's expression text with a new one, in\n// left-to-right, top-to-bottom (row-major) order. Only the text\n// content changes -- run/paragraph formatting (font, size, alignment)\n// must stay untouched.\n\nconst newValues = [\n  \"54-29=\", \"12-9=\", \"30+25=\", \"85+14=\", \"83-67=\",\n  \"42+36=\", \"90-42=\", \"47-43=\", \"72-54=\", \"72-8=\",\n  \"45+48=\", \"52+7=\", \"9+21=\", \"38+30=\", \"8+32=\",\n  \"8+34=\", \"83+16=\", \"31-20=\", \"27-5=\", \"17+79=\",\n  \"95-3=\", \"84+11=\", \"75-16=\", \"70+23=\", \"17+78=\",\n  \"28+38=\", \"50+35=\", \"5+65=\", \"24+54=\", \"44+7=\",\n  \"98-71=\", \"25+34=\", \"92-1=\", \"94-86=\", \"89-2=\",\n  \"62+17=\", \"82-35=\", \"86+1=\", \"56+17=\", \"97-67=\",\n  \"42-39=\", \"20+50=\", \"32-12=\", \"71-7=\", \"33+17=\",\n  \"9+66=\", \"42+36=\", \"26+63=\", \"2+92=\", \"10+69=\",\n  \"47+1=\", \"33+60=\", \"30+34=\", \"0+28=\", \"52+46=\",\n  \"55+6=\", \"15+61=\", \"83-8=\", \"6+11=\", \"33-12=\",\n  \"7+9=\", \"79-64=\", \"13+74=\", \"95-86=\", \"57-38=\",\n  \"89-1=\", \"47+51=\", \"14+76=\", \"69-69=\", \"27+46=\",\n  \"99-86=\", \"3+68=\", \"59+15=\", \"63+27=\", \"51+0=\",\n  \"24+37=\", \"63+8=\", \"68-21=\", \"11-0=\", \"74-48=\",\n  \"77-10=\", \"67+14=\", \"24+63=\", \"99-93=\", \"56+2=\",\n  \"7+45=\", \"41-15=\", \"12+50=\", \"87-43=\", \"22+67=\",\n  \"35+49=\", \"35+51=\", \"70-69=\", \"50-39=\", \"74-37=\",\n  \"36+37=\", \"3+35=\", \"77-42=\", \"61-28=\", \"67-41=\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet i = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (const cell of cells.items) {\n    if (i >= newValues.length) break;\n    cell.value = newValues[i];\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of arithmetic expressions\n# (20 rows x 5 columns = 100 cells, one run/paragraph per cell). The\n# commit replaces each cell's expression text with a new one, in\n# left-to-right, top-to-bottom (row-major) order. Only the text\n# content changes -- run/paragraph formatting (font, size, alignment)\n# must stay untouched, so we assign straight to Cell.Range.Text\n# (Word COM strips/re-adds the trailing cell mark automatically).\n\n$newValues = @(\n  \"54-29=\", \"12-9=\", \"30+25=\", \"85+14=\", \"83-67=\",\n  \"42+36=\", \"90-42=\", \"47-43=\", \"72-54=\", \"72-8=\",\n  \"45+48=\", \"52+7=\", \"9+21=\", \"38+30=\", \"8+32=\",\n  \"8+34=\", \"83+16=\", \"31-20=\", \"27-5=\", \"17+79=\",\n  \"95-3=\", \"84+11=\", \"75-16=\", \"70+23=\", \"17+78=\",\n  \"28+38=\", \"50+35=\", \"5+65=\", \"24+54=\", \"44+7=\",\n  \"98-71=\", \"25+34=\", \"92-1=\", \"94-86=\", \"89-2=\",\n  \"62+17=\", \"82-35=\", \"86+1=\", \"56+17=\", \"97-67=\",\n  \"42-39=\", \"20+50=\", \"32-12=\", \"71-7=\", \"33+17=\",\n  \"9+66=\", \"42+36=\", \"26+63=\", \"2+92=\", \"10+69=\",\n  \"47+1=\", \"33+60=\", \"30+34=\", \"0+28=\", \"52+46=\",\n  \"55+6=\", \"15+61=\", \"83-8=\", \"6+11=\", \"33-12=\",\n  \"7+9=\", \"79-64=\", \"13+74=\", \"95-86=\", \"57-38=\",\n  \"89-1=\", \"47+51=\", \"14+76=\", \"69-69=\", \"27+46=\",\n  \"99-86=\", \"3+68=\", \"59+15=\", \"63+27=\", \"51+0=\",\n  \"24+37=\", \"63+8=\", \"68-21=\", \"11-0=\", \"74-48=\",\n  \"77-10=\", \"67+14=\", \"24+63=\", \"99-93=\", \"56+2=\",\n  \"7+45=\", \"41-15=\", \"12+50=\", \"87-43=\", \"22+67=\",\n  \"35+49=\", \"35+51=\", \"70-69=\", \"50-39=\", \"74-37=\",\n  \"36+37=\", \"3+35=\", \"77-42=\", \"61-28=\", \"67-41=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$numRows = $t.Rows.Count\n$numCols = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $numRows; $r++) {\n  for ($c = 1; $c -le $numCols; $c++) {\n    if ($i -ge $newValues.Length) { break }\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$i]\n    $i++\n  }\n}\n"}
